{"js": "// Replace the date line and every \"NNxNN=NNNN\" multiplication answer in the\n// table with its updated value. All 26 text runs in the document change, and\n// none of the \"from\" values repeat, so a straightforward search/replace per\n// pair is unambiguous.\nconst pairs = [\n  [\"2025-09-24 Wednesday\", \"2025-09-25 Thursday\"],\n  [\"70\u00d768=4760\", \"71\u00d769=4899\"],\n  [\"21\u00d743=903\", \"31\u00d777=2387\"],\n  [\"36\u00d743=1548\", \"37\u00d759=2183\"],\n  [\"48\u00d761=2928\", \"86\u00d799=8514\"],\n  [\"25\u00d746=1150\", \"37\u00d792=3404\"],\n  [\"88\u00d747=4136\", \"38\u00d755=2090\"],\n  [\"32\u00d784=2688\", \"76\u00d738=2888\"],\n  [\"96\u00d774=7104\", \"84\u00d799=8316\"],\n  [\"27\u00d749=1323\", \"84\u00d740=3360\"],\n  [\"48\u00d740=1920\", \"55\u00d755=3025\"],\n  [\"93\u00d788=8184\", \"29\u00d712=348\"],\n  [\"68\u00d760=4080\", \"87\u00d780=6960\"],\n  [\"87\u00d762=5394\", \"37\u00d741=1517\"],\n  [\"22\u00d743=946\", \"78\u00d713=1014\"],\n  [\"74\u00d759=4366\", \"54\u00d764=3456\"],\n  [\"82\u00d752=4264\", \"81\u00d779=6399\"],\n  [\"88\u00d735=3080\", \"85\u00d766=5610\"],\n  [\"94\u00d780=7520\", \"57\u00d755=3135\"],\n  [\"39\u00d725=975\", \"76\u00d771=5396\"],\n  [\"26\u00d795=2470\", \"57\u00d751=2907\"],\n  [\"46\u00d752=2392\", \"31\u00d777=2387\"],\n  [\"66\u00d716=1056\", \"61\u00d726=1586\"],\n  [\"11\u00d714=154\", \"43\u00d772=3096\"],\n  [\"35\u00d746=1610\", \"79\u00d751=4029\"],\n  [\"79\u00d766=5214\", \"88\u00d761=5368\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of pairs) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${from}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"NNxNN=NNNN\" multiplication answer in the\n# table with its updated value. All 26 text runs in the document change, and\n# none of the \"from\" values repeat, so a straightforward Find/Replace per\n# pair is unambiguous.\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$pairs = @(\n    ,@(\"2025-09-24 Wednesday\", \"2025-09-25 Thursday\")\n    ,@(\"70\u00d768=4760\", \"71\u00d769=4899\")\n    ,@(\"21\u00d743=903\", \"31\u00d777=2387\")\n    ,@(\"36\u00d743=1548\", \"37\u00d759=2183\")\n    ,@(\"48\u00d761=2928\", \"86\u00d799=8514\")\n    ,@(\"25\u00d746=1150\", \"37\u00d792=3404\")\n    ,@(\"88\u00d747=4136\", \"38\u00d755=2090\")\n    ,@(\"32\u00d784=2688\", \"76\u00d738=2888\")\n    ,@(\"96\u00d774=7104\", \"84\u00d799=8316\")\n    ,@(\"27\u00d749=1323\", \"84\u00d740=3360\")\n    ,@(\"48\u00d740=1920\", \"55\u00d755=3025\")\n    ,@(\"93\u00d788=8184\", \"29\u00d712=348\")\n    ,@(\"68\u00d760=4080\", \"87\u00d780=6960\")\n    ,@(\"87\u00d762=5394\", \"37\u00d741=1517\")\n    ,@(\"22\u00d743=946\", \"78\u00d713=1014\")\n    ,@(\"74\u00d759=4366\", \"54\u00d764=3456\")\n    ,@(\"82\u00d752=4264\", \"81\u00d779=6399\")\n    ,@(\"88\u00d735=3080\", \"85\u00d766=5610\")\n    ,@(\"94\u00d780=7520\", \"57\u00d755=3135\")\n    ,@(\"39\u00d725=975\", \"76\u00d771=5396\")\n    ,@(\"26\u00d795=2470\", \"57\u00d751=2907\")\n    ,@(\"46\u00d752=2392\", \"31\u00d777=2387\")\n    ,@(\"66\u00d716=1056\", \"61\u00d726=1586\")\n    ,@(\"11\u00d714=154\", \"43\u00d772=3096\")\n    ,@(\"35\u00d746=1610\", \"79\u00d751=4029\")\n    ,@(\"79\u00d766=5214\", \"88\u00d761=5368\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $from = $pair[0]\n    $to = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $from\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $to\n    $ok = $find.Execute($null, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $null, $wdReplaceAll)\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $from\"\n    }\n}\n"}
